$p = $ppt.ActivePresentation

$newStyleId = "{18CB0E6E-7971-41E5-8EA3-7CC97A1CB6F4}"

foreach ($slideIdx in 14, 15, 16) {
    $s = $p.Slides.Item($slideIdx)
    $shp = $s.Shapes.Item(1)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle($newStyleId)
    }
}
